$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row data for alloy entries (Al-Si LPBF series, sciencedirect S1359645421000781) ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "https://www.sciencedirect.com/science/article/pii/S1359645421000781"
$ws.Range("C2").Value = "5a"
$ws.Range("D2").Value = "Al-0.5Si"
$ws.Range("E2").Value = "Aluminum"
$ws.Range("F2").Value = 99.5
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0.5
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AP2").Value = 0
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = 0
$ws.Range("AS2").Value = 0
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = "LPBF"
$ws.Range("AW2").Value = 350
$ws.Range("AX2").Value = 1200
$ws.Range("AY2").Value = 130
$ws.Range("AZ2").Value = 30
$ws.Range("BA2").Value = 74.7863247863248

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "https://www.sciencedirect.com/science/article/pii/S1359645421000781"
$ws.Range("C3").Value = "5a"
$ws.Range("D3").Value = "Al-1.0Si"
$ws.Range("E3").Value = "Aluminum"
$ws.Range("F3").Value = 99
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 0
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 1
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 0
$ws.Range("AK3").Value = 0
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AP3").Value = 1
$ws.Range("AQ3").Value = 101.248
$ws.Range("AR3").Value = 40.36709
$ws.Range("AS3").Value = 11
$ws.Range("AT3").Value = 0.14814
$ws.Range("AU3").Value = 2.9974212906709865
$ws.Range("AV3").Value = "LPBF"
$ws.Range("AW3").Value = 350
$ws.Range("AX3").Value = 1200
$ws.Range("AY3").Value = 130
$ws.Range("AZ3").Value = 30
$ws.Range("BA3").Value = 74.7863247863248

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "https://www.sciencedirect.com/science/article/pii/S1359645421000781"
$ws.Range("C4").Value = "5a"
$ws.Range("D4").Value = "Al-2.0Si"
$ws.Range("E4").Value = "Aluminum"
$ws.Range("F4").Value = 98
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 2
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AK4").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AP4").Value = 0
$ws.Range("AQ4").Value = 0
$ws.Range("AR4").Value = 0
$ws.Range("AS4").Value = 0
$ws.Range("AU4").Value = 0
$ws.Range("AV4").Value = "LPBF"
$ws.Range("AW4").Value = 350
$ws.Range("AX4").Value = 1200
$ws.Range("AY4").Value = 130
$ws.Range("AZ4").Value = 30
$ws.Range("BA4").Value = 74.7863247863248

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "https://www.sciencedirect.com/science/article/pii/S1359645421000781"
$ws.Range("C5").Value = "5a"
$ws.Range("D5").Value = "Al-4.0Si"
$ws.Range("E5").Value = "Aluminum"
$ws.Range("F5").Value = 96
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 4
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 0
$ws.Range("AK5").Value = 0
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AP5").Value = 0
$ws.Range("AQ5").Value = 0
$ws.Range("AR5").Value = 0
$ws.Range("AS5").Value = 0
$ws.Range("AU5").Value = 0
$ws.Range("AV5").Value = "LPBF"
$ws.Range("AW5").Value = 350
$ws.Range("AX5").Value = 1200
$ws.Range("AY5").Value = 130
$ws.Range("AZ5").Value = 30
$ws.Range("BA5").Value = 74.7863247863248

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "https://www.sciencedirect.com/science/article/pii/S1359645421000781"
$ws.Range("C6").Value = "5a"
$ws.Range("D6").Value = "Al-12.6Si"
$ws.Range("E6").Value = "Aluminum"
$ws.Range("F6").Value = 87.4
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 12.6
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 0
$ws.Range("AK6").Value = 0
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AP6").Value = 0
$ws.Range("AQ6").Value = 0
$ws.Range("AR6").Value = 0
$ws.Range("AS6").Value = 0
$ws.Range("AU6").Value = 0
$ws.Range("AV6").Value = "LPBF"
$ws.Range("AW6").Value = 350
$ws.Range("AX6").Value = 1200
$ws.Range("AY6").Value = 130
$ws.Range("AZ6").Value = 30
$ws.Range("BA6").Value = 74.7863247863248

# Row 7
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "https://www.sciencedirect.com/science/article/pii/S1359645421000781"
$ws.Range("C7").Value = "5a"
$ws.Range("D7").Value = "Al-16.0Si"
$ws.Range("E7").Value = "Aluminum"
$ws.Range("F7").Value = 84
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0
$ws.Range("AC7").Value = 0
$ws.Range("AD7").Value = 16
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = 0
$ws.Range("AK7").Value = 0
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AP7").Value = 0
$ws.Range("AQ7").Value = 0
$ws.Range("AR7").Value = 0
$ws.Range("AS7").Value = 0
$ws.Range("AU7").Value = 0
$ws.Range("AV7").Value = "LPBF"
$ws.Range("AW7").Value = 350
$ws.Range("AX7").Value = 1200
$ws.Range("AY7").Value = 130
$ws.Range("AZ7").Value = 30
$ws.Range("BA7").Value = 74.7863247863248

# --- Row 2 style normalization to match the rest of the templated rows ---
# (AS2 picks up the integer-count format, BB2 picks up the left-aligned blank style)
$ws.Range("AS2").NumberFormat = "0"
$ws.Range("BB2").HorizontalAlignment = -4131

# --- Hyperlink the Reference cell (B2) to the source article ---
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.sciencedirect.com/science/article/pii/S1359645421000781")

# --- Restore selection to the cell the author left active ---
$ws.Range("L13").Select()
